$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the attendance detection data that was recorded for row 4
# (Rahul Jaluthria), restoring it to the blank/unmarked state shared
# by the other students' rows.
$ws.Range("D4:L4").ClearContents()
